$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.166.25"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "'2.998.58"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'541.88"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").Value = "'139.86"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.00%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'2.996.60"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").Value = "'0.488"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D10").Value = "'6.78"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +12.28%  "
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").Value = "'0.444"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").Value = "'0.0000219"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "'33.95"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").Value = "'3.483.50"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "'62.320.35"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "'3.005.64"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'0.107"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.15%  "
$ws.Range("D19").Value = "'6.54"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("D20").Value = "'466.42"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("D21").Value = "'13.37"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").Value = "'0.651"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.35%  "
$ws.Range("D23").Value = "'7.18"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").Value = "'79.37"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("D25").Value = "'12.58"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.99%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "'2.71"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("D28").Value = "'7.60"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.01%  "
$ws.Range("D29").Value = "'2.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.92%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").Value = "'25.41"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("E32").Value = "  -2.31%  "
$ws.Range("D33").Value = "'2.33"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("D34").Value = "'5.54"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").Value = "'54.68"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("D36").Value = "'5.81"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("D37").Value = "'449.17"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.73%  "
$ws.Range("D38").Value = "'0.0807"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("D39").Value = "'0.0388"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("D40").Value = "'2.944.54"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -8.46%  "
$ws.Range("E41").Value = "  -3.68%  "
$ws.Range("D42").Value = "'8.05"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("D43").Value = "'2.57"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.92%  "
$ws.Range("D44").Value = "'26.72"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D46").Value = "'0.246"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("D47").Value = "'2.00"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").Value = "'115.13"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("D50").Value = "'0.0₃0494"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("B51").Value = "BitgetToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D51").Value = "'1.24"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.30%  "
